$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A260").Value = 'Religious organizations'
$ws.Range("B260").Value = 'Fees and admissions'

$ws.Range("A261").Value = 'Legal services'
$ws.Range("B261").Value = 'Other household expenses'

$ws.Range("A262").Value = 'Accounting, tax preparation, bookkeeping, and payroll services'
$ws.Range("B262").Value = 'Other household expenses'

$ws.Range("A263").Value = 'Specialized design services'
$ws.Range("B263").Value = 'Other household expenses'

$ws.Range("A264").Value = 'Scientific research and development services'
$ws.Range("B264").Value = 'Other household expenses'

$ws.Range("A265").Value = 'Advertising, public relations, and related services'
$ws.Range("B265").Value = 'Other household expenses'

$ws.Range("A266").Value = 'Funds, trusts, and other financial vehicles'
$ws.Range("B266").Value = 'Other household expenses'

$ws.Range("A267").Value = 'Securities and commodity contracts intermediation and brokerage'
$ws.Range("B267").Value = 'Other household expenses'

$ws.Range("A268").Value = 'Other financial investment activities'
$ws.Range("B268").Value = 'Other household expenses'

$ws.Range("A269").Value = 'Couriers and messengers'
$ws.Range("B269").Value = 'Other household expenses'

$ws.Range("A270").Value = 'Noncomparable imports'
$ws.Range("B270").Value = 'Other household expenses'

$ws.Range("A271").Value = 'Private households'
$ws.Range("B271").Value = 'Other household expenses'

$ws.Range("A272").Value = 'Other state and local government enterprises'
$ws.Range("B272").Value = 'Other household expenses'

$ws.Range("A273").Value = 'Individual and family services'
$ws.Range("B273").Value = 'Other household expenses'

$ws.Range("A274").Value = 'Other support services'
$ws.Range("B274").Value = 'Other household expenses'

$ws.Range("A275").Value = 'Veterinary services'
$ws.Range("B275").Value = 'Other household expenses'

$ws.Range("A276").Value = 'Employment services'
$ws.Range("B276").Value = 'Other household expenses'

$ws.Range("A277").Value = 'Business support services'
$ws.Range("B277").Value = 'Other household expenses'

$ws.Range("A278").Value = 'Travel arrangement and reservation services'
$ws.Range("B278").Value = 'Other household expenses'

$ws.Range("A279").Value = 'Investigation and security services'
$ws.Range("B279").Value = 'Other household expenses'

$ws.Range("A280").Value = 'Other real estate'
$ws.Range("B280").Value = 'Mortgage interest and charges'

$ws.Range("A281").Value = 'Pipeline transportation'
$ws.Range("B281").Value = 'Fuel oil and other fuels'

$ws.Range("A282").Value = 'Support activities for agriculture and forestry'
$ws.Range("B282").Value = 'Water and other public services'

$ws.Range("A283").Value = 'Grantmaking, giving, and social advocacy organizations'
$ws.Range("B283").Value = 'Other household expenses'

$ws.Range("A284").Value = 'Civic, social, professional, and similar organizations'
$ws.Range("B284").Value = 'Other household expenses'

$ws.Range("A285").Value = 'Commercial and industrial machinery and equipment rental and leasing'
$ws.Range("B285").Value = 'Other household expenses'

$ws.Range("A286").Value = 'Warehousing and storage'
$ws.Range("B286").Value = 'Other household expenses'

$ws.Range("A287").Value = 'Wholesale trade'
$ws.Range("B287").Value = 'Other household expenses'

$ws.Range("A288").Value = 'Other nonmetallic mineral mining and quarrying'
$ws.Range("B288").Value = 'Major appliances'
$ws.Range("C288").Value = 'Small appliances, miscellaneous housewares'

$ws.Range("A289").Value = 'Iron and steel mills and ferroalloy manufacturing'
$ws.Range("B289").Value = 'Major appliances'
$ws.Range("C289").Value = 'Small appliances, miscellaneous housewares'

$ws.Range("A290").Value = 'Nonferrous metal (except copper and aluminum) rolling, drawing, extruding and alloying'
$ws.Range("B290").Value = 'Major appliances'
$ws.Range("C290").Value = 'Small appliances, miscellaneous housewares'

$ws.Range("A291").Value = 'Nonferrous metal foundries'
$ws.Range("B291").Value = 'Major appliances'
$ws.Range("C291").Value = 'Small appliances, miscellaneous housewares'

$ws.Range("A292").Value = 'Crown and closure manufacturing and metal stamping'
$ws.Range("B292").Value = 'Major appliances'
$ws.Range("C292").Value = 'Small appliances, miscellaneous housewares'

$ws.Range("A293").Value = 'Plate work and fabricated structural product manufacturing'
$ws.Range("B293").Value = 'Major appliances'
$ws.Range("C293").Value = 'Small appliances, miscellaneous housewares'

$ws.Range("A294").Value = 'Metal can, box, and other metal container (light gauge) manufacturing'
$ws.Range("B294").Value = 'Major appliances'
$ws.Range("C294").Value = 'Small appliances, miscellaneous housewares'

$ws.Range("A295").Value = 'Hardware manufacturing'
$ws.Range("B295").Value = 'Major appliances'
$ws.Range("C295").Value = 'Small appliances, miscellaneous housewares'

$ws.Range("A296").Value = 'Spring and wire product manufacturing'
$ws.Range("B296").Value = 'Major appliances'
$ws.Range("C296").Value = 'Small appliances, miscellaneous housewares'

$ws.Range("A297").Value = 'Office machinery manufacturing'
$ws.Range("B297").Value = 'Major appliances'
$ws.Range("C297").Value = 'Small appliances, miscellaneous housewares'

$ws.Range("A298").Value = 'Metal cutting and forming machine tool manufacturing'
$ws.Range("B298").Value = 'Major appliances'
$ws.Range("C298").Value = 'Small appliances, miscellaneous housewares'

$ws.Range("A299").Value = 'Other engine equipment manufacturing'
$ws.Range("B299").Value = 'Major appliances'
$ws.Range("C299").Value = 'Small appliances, miscellaneous housewares'

$ws.Range("A300").Value = 'Industrial gas manufacturing'
$ws.Range("B300").Value = 'Major appliances'
$ws.Range("C300").Value = 'Small appliances, miscellaneous housewares'

$ws.Range("C300").Select()